$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.010371260139372
$ws.Range("D2").Value = 1.012902601556871
$ws.Range("E2").Value = 1.012676467692917
$ws.Range("F2").Value = 1.022948532896612
$ws.Range("I2").Value = 1.025565289030852
$ws.Range("J2").Value = 1.015625892347066
$ws.Range("K2").Value = 1.015765362915779
$ws.Range("L2").Value = 1.015539909156081
$ws.Range("M2").Value = 1.025781411076117
$ws.Range("N2").Value = 1.009309119679293
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.011406944684967
$ws.Range("D3").Value = 1.013798066038265
$ws.Range("E3").Value = 1.013557327623777
$ws.Range("F3").Value = 1.02415355767728
$ws.Range("I3").Value = 1.02553972418465
$ws.Range("J3").Value = 1.016293901614334
$ws.Range("K3").Value = 1.01646518686554
$ws.Range("L3").Value = 1.016225119791121
$ws.Range("M3").Value = 1.02679212406856
$ws.Range("N3").Value = 1.009535773716937
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.012077623653653
$ws.Range("D4").Value = 1.014378277196461
$ws.Range("E4").Value = 1.01412811510856
$ws.Range("F4").Value = 1.0249323242737
$ws.Range("I4").Value = 1.02552073760344
$ws.Range("J4").Value = 1.016726148005834
$ws.Range("K4").Value = 1.016918179526824
$ws.Range("L4").Value = 1.016668679712447
$ws.Range("M4").Value = 1.027444604017036
$ws.Range("N4").Value = 1.009682273690448
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.012359702064887
$ws.Range("D5").Value = 1.014622385505205
$ws.Range("E5").Value = 1.014368267983479
$ws.Range("F5").Value = 1.025259486731346
$ws.Range("I5").Value = 1.025512169486292
$ws.Range("J5").Value = 1.016907864006635
$ws.Range("K5").Value = 1.017108655500331
$ws.Range("L5").Value = 1.016855195648634
$ws.Range("M5").Value = 1.027718542588932
$ws.Range("N5").Value = 1.009743823679129
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.012407071601336
$ws.Range("D6").Value = 1.014663383344966
$ws.Range("E6").Value = 1.014408602055124
$ws.Range("F6").Value = 1.025314405200366
$ws.Range("I6").Value = 1.025510696460925
$ws.Range("J6").Value = 1.016938374891133
$ws.Range("K6").Value = 1.017140639463113
$ws.Range("L6").Value = 1.016886515024525
$ws.Range("M6").Value = 1.027764516724678
$ws.Range("N6").Value = 1.009754155923315
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.012081392308069
$ws.Range("D7").Value = 1.014381538247465
$ws.Range("E7").Value = 1.01413132328226
$ws.Range("F7").Value = 1.024936696742331
$ws.Range("I7").Value = 1.02552062542028
$ws.Range("J7").Value = 1.016728576105841
$ws.Range("K7").Value = 1.016920724527846
$ws.Range("L7").Value = 1.016671171777739
$ws.Range("M7").Value = 1.02744826582848
$ws.Range("N7").Value = 1.009683096276352
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.01072116656477
$ws.Range("D8").Value = 1.013205064928666
$ws.Range("E8").Value = 1.012973989880497
$ws.Range("F8").Value = 1.023355976484694
$ws.Range("I8").Value = 1.025557154830652
$ws.Range("J8").Value = 1.015851649072877
$ws.Range("K8").Value = 1.016001838380274
$ws.Range("L8").Value = 1.015771441190121
$ws.Range("M8").Value = 1.026123300559184
$ws.Range("N8").Value = 1.009385751544963
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.008328268902935
$ws.Range("D9").Value = 1.01113801167686
$ws.Range("E9").Value = 1.010940868525497
$ws.Range("F9").Value = 1.020563152713118
$ws.Range("I9").Value = 1.025602844781703
$ws.Range("J9").Value = 1.014306403174602
$ws.Range("K9").Value = 1.014383884164755
$ws.Range("L9").Value = 1.014187419502133
$ws.Range("M9").Value = 1.023776923686797
$ws.Range("N9").Value = 1.008860573384166
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.006735682695536
$ws.Range("D10").Value = 1.00976407107002
$ws.Range("E10").Value = 1.009589690254689
$ws.Range("F10").Value = 1.018696291817971
$ws.Range("I10").Value = 1.025620801168555
$ws.Range("J10").Value = 1.013276263651415
$ws.Range("K10").Value = 1.013306103196133
$ws.Range("L10").Value = 1.013132381908735
$ws.Range("M10").Value = 1.022204875520217
$ws.Range("N10").Value = 1.008509647128149
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.006046708396777
$ws.Range("D11").Value = 1.00917011698815
$ws.Range("E11").Value = 1.009005626708028
$ws.Range("F11").Value = 1.017886736216683
$ws.Range("I11").Value = 1.025625622239129
$ws.Range("J11").Value = 1.012830210575164
$ws.Range("K11").Value = 1.012839618851596
$ws.Range("L11").Value = 1.012675775016912
$ws.Range("M11").Value = 1.021522313630841
$ws.Range("N11").Value = 1.008357502984775
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.005790886424283
$ws.Range("D12").Value = 1.008949642345352
$ws.Range("E12").Value = 1.008788831066341
$ws.Range("F12").Value = 1.017585851580468
$ws.Range("I12").Value = 1.025626970047317
$ws.Range("J12").Value = 1.012664527388108
$ws.Range("K12").Value = 1.012666376378642
$ws.Range("L12").Value = 1.012506205930539
$ws.Range("M12").Value = 1.021268501066274
$ws.Range("N12").Value = 1.008300961378916
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.005845756867743
$ws.Range("D13").Value = 1.008996928275204
$ws.Range("E13").Value = 1.008835327600911
$ws.Range("F13").Value = 1.01765040049325
$ws.Range("I13").Value = 1.025626700968447
$ws.Range("J13").Value = 1.01270006695595
$ws.Range("K13").Value = 1.012703536089303
$ws.Range("L13").Value = 1.012542577484794
$ws.Range("M13").Value = 1.021322957335061
$ws.Range("N13").Value = 1.008313091032246
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.006025560161781
$ws.Range("D14").Value = 1.009151889497459
$ws.Range("E14").Value = 1.008987703217944
$ws.Range("F14").Value = 1.017861868666231
$ws.Range("I14").Value = 1.025625742680436
$ws.Range("J14").Value = 1.012816515134254
$ws.Range("K14").Value = 1.012825297949441
$ws.Range("L14").Value = 1.012661757666803
$ws.Range("M14").Value = 1.021501339090553
$ws.Range("N14").Value = 1.008352829818385
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.006136355320617
$ws.Range("D15").Value = 1.009247385626276
$ws.Range("E15").Value = 1.009081606963117
$ws.Range("F15").Value = 1.017992137368239
$ws.Range("I15").Value = 1.025625093578135
$ws.Range("J15").Value = 1.012888262809247
$ws.Range("K15").Value = 1.012900323509974
$ws.Range("L15").Value = 1.012735193162555
$ws.Range("M15").Value = 1.021611209046949
$ws.Range("N15").Value = 1.008377310424818
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.006781420815183
$ws.Range("D16").Value = 1.00980351034626
$ws.Range("E16").Value = 1.009628473856651
$ws.Range("F16").Value = 1.018749994191611
$ws.Range("I16").Value = 1.025620419044362
$ws.Range("J16").Value = 1.013305866826242
$ws.Range("K16").Value = 1.013337066501099
$ws.Range("L16").Value = 1.013162690303038
$ws.Range("M16").Value = 1.022250135825731
$ws.Range("N16").Value = 1.008519740428696
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.007186220755677
$ws.Range("D17").Value = 1.010152613028368
$ws.Range("E17").Value = 1.009971778575196
$ws.Range("F17").Value = 1.019225058152309
$ws.Range("I17").Value = 1.025616696543157
$ws.Range("J17").Value = 1.013567820016914
$ws.Range("K17").Value = 1.013611078076223
$ws.Range("L17").Value = 1.013430910105915
$ws.Range("M17").Value = 1.022650421141879
$ws.Range("N17").Value = 1.008609032050705
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00742239413564
$ws.Range("D18").Value = 1.01035633240114
$ws.Range("E18").Value = 1.010172119398475
$ws.Range("F18").Value = 1.019502039895548
$ws.Range("I18").Value = 1.025614240179127
$ws.Range("J18").Value = 1.013720613269095
$ws.Range("K18").Value = 1.013770923858688
$ws.Range("L18").Value = 1.013587380481613
$ws.Range("M18").Value = 1.022883721850529
$ws.Range("N18").Value = 1.008661095926818
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.00750293341286
$ws.Range("D19").Value = 1.01042581132021
$ws.Range("E19").Value = 1.010240446878046
$ws.Range("F19").Value = 1.019596464015865
$ws.Range("I19").Value = 1.025613354239018
$ws.Range("J19").Value = 1.013772711880347
$ws.Range("K19").Value = 1.013825430434121
$ws.Range("L19").Value = 1.013640736610361
$ws.Range("M19").Value = 1.022963241012383
$ws.Range("N19").Value = 1.008678845234438
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.007142783305193
$ws.Range("D20").Value = 1.010115147921038
$ws.Range("E20").Value = 1.009934935192225
$ws.Range("F20").Value = 1.01917410018625
$ws.Range("I20").Value = 1.025617125414324
$ws.Range("J20").Value = 1.01353971487725
$ws.Range("K20").Value = 1.013581677212053
$ws.Range("L20").Value = 1.013402130350465
$ws.Range("M20").Value = 1.022607492824765
$ws.Range("N20").Value = 1.008599453808395
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.005972609987333
$ws.Range("D21").Value = 1.009106253236438
$ws.Range("E21").Value = 1.008942828208719
$ws.Range("F21").Value = 1.017799601522866
$ws.Range("I21").Value = 1.025626037093034
$ws.Range("J21").Value = 1.012782224025121
$ws.Range("K21").Value = 1.012789441270261
$ws.Range("L21").Value = 1.012626661104018
$ws.Range("M21").Value = 1.021448817803202
$ws.Range("N21").Value = 1.00834112852168
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.005237417892617
$ws.Range("D22").Value = 1.008472768046197
$ws.Range("E22").Value = 1.00831992832184
$ws.Range("F22").Value = 1.016934360067156
$ws.Range("I22").Value = 1.025629077859853
$ws.Range("J22").Value = 1.012305964090164
$ws.Range("K22").Value = 1.012291508153685
$ws.Range("L22").Value = 1.012139295703456
$ws.Range("M22").Value = 1.020718700393434
$ws.Range("N22").Value = 1.008178544367256
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.005627105825429
$ws.Range("D23").Value = 1.008808510109015
$ws.Range("E23").Value = 1.008650056030312
$ws.Range("F23").Value = 1.017393139639279
$ws.Range("I23").Value = 1.025627708484345
$ws.Range("J23").Value = 1.012558437999322
$ws.Range("K23").Value = 1.012555455054654
$ws.Range("L23").Value = 1.012397637996348
$ws.Range("M23").Value = 1.021105902306827
$ws.Range("N23").Value = 1.008264748870206
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.007162410632428
$ws.Range("D24").Value = 1.010132076500852
$ws.Range("E24").Value = 1.009951582831367
$ws.Range("F24").Value = 1.019197126252603
$ws.Range("I24").Value = 1.025616932506859
$ws.Range("J24").Value = 1.013552414378807
$ws.Range("K24").Value = 1.013594962136805
$ws.Range("L24").Value = 1.01341513461417
$ws.Range("M24").Value = 1.022626890835974
$ws.Range("N24").Value = 1.008603781860695
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.008946417674851
$ws.Range("D25").Value = 1.011671675065826
$ws.Range("E25").Value = 1.011465735081954
$ws.Range("F25").Value = 1.021286041431387
$ws.Range("I25").Value = 1.02559324068728
$ws.Range("J25").Value = 1.014705883189532
$ws.Range("K25").Value = 1.014802014920101
$ws.Range("L25").Value = 1.014596756585115
$ws.Range("M25").Value = 1.024384893637066
$ws.Range("N25").Value = 1.008996487722042
